$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# "Tabelle1" -> "Früchte"
$ws1.Name = "Früchte"

# Update the data value and the (now inactive) selection on the first sheet.
$ws1.Range("C9").Value = 31
$ws1.Range("C2").Select()

# "Tabelle2" -> "Tabelle1", and make it the active sheet/tab.
$ws2.Name = "Tabelle1"
$ws2.Activate()
